$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, [string]$value) {
    # Leading apostrophe forces Excel to store the value as literal text
    # (prevents auto-conversion of numeric-looking strings like "0.999").
    $cell.Value = "'" + $value
    # Re-apply the default style so the quote-prefix flag added above
    # doesn't leave a stray number-format / style behind.
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "61.573.49"
Set-TextCell $ws.Range("E2") "  +0.98%  "

# Row 3
Set-TextCell $ws.Range("D3") "3.390.91"
Set-TextCell $ws.Range("E3") "  -0.26%  "

# Row 4
Set-TextCell $ws.Range("D4") "0.999"
Set-TextCell $ws.Range("E4") "  -0.01%  "

# Row 5
Set-TextCell $ws.Range("D5") "575.76"
Set-TextCell $ws.Range("E5") "  +0.46%  "

# Row 6
Set-TextCell $ws.Range("D6") "141.25"
Set-TextCell $ws.Range("E6") "  -0.97%  "

# Row 7
Set-TextCell $ws.Range("E7") "  +0.02%  "

# Row 8
Set-TextCell $ws.Range("E8") "  -0.74%  "

# Row 9
Set-TextCell $ws.Range("D9") "7.67"
Set-TextCell $ws.Range("E9") "  +1.23%  "

# Row 10
Set-TextCell $ws.Range("E10") "  -1.23%  "

# Row 11
Set-TextCell $ws.Range("E11") "  -2.52%  "

# Row 12
Set-TextCell $ws.Range("D12") "3.969.35"
Set-TextCell $ws.Range("E12") "  -0.22%  "

# Row 13
Set-TextCell $ws.Range("E13") "  +0.17%  "

# Row 14
Set-TextCell $ws.Range("D14") "28.41"
Set-TextCell $ws.Range("E14") "  +1.14%  "

# Row 15
Set-TextCell $ws.Range("D15") "3.382.17"
Set-TextCell $ws.Range("E15") "  -0.18%  "

# Row 16
Set-TextCell $ws.Range("E16") "  -0.73%  "

# Row 17
Set-TextCell $ws.Range("D17") "61.618.93"
Set-TextCell $ws.Range("E17") "  +0.95%  "

# Row 18
Set-TextCell $ws.Range("E18") "  -0.18%  "

# Row 19
Set-TextCell $ws.Range("D19") "13.62"
Set-TextCell $ws.Range("E19") "  -1.83%  "

# Row 20
Set-TextCell $ws.Range("D20") "9.00"
Set-TextCell $ws.Range("E20") "  +0.32%  "

# Row 21
Set-TextCell $ws.Range("D21") "391.64"
Set-TextCell $ws.Range("E21") "  +2.14%  "

# Row 22
Set-TextCell $ws.Range("E22") "  +0.73%  "

# Row 23
Set-TextCell $ws.Range("E23") "  -1.64%  "

# Row 24
Set-TextCell $ws.Range("E24") "  +0.01%  "

# Row 25
Set-TextCell $ws.Range("D25") "0.0000113"
Set-TextCell $ws.Range("E25") "  -4.38%  "

# Row 26
Set-TextCell $ws.Range("E26") "  +7.80%  "

# Row 27
Set-TextCell $ws.Range("E27") "  -0.05%  "

# Row 28
Set-TextCell $ws.Range("D28") "7.27"
Set-TextCell $ws.Range("E28") "  -1.48%  "

# Row 29
Set-TextCell $ws.Range("E29") "  +0.02%  "

# Row 30
Set-TextCell $ws.Range("E30") "  -0.94%  "

# Row 31
Set-TextCell $ws.Range("E31") "  -1.29%  "

# Row 32
Set-TextCell $ws.Range("E32") "  +0.08%  "

# Row 33
Set-TextCell $ws.Range("D33") "23.30"
Set-TextCell $ws.Range("E33") "  -0.85%  "

# Row 34
Set-TextCell $ws.Range("B34") "Monero"
Set-TextCell $ws.Range("C34") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws.Range("D34") "169.48"
Set-TextCell $ws.Range("E34") "  +1.04%  "

# Row 35
Set-TextCell $ws.Range("B35") "Aptos"
Set-TextCell $ws.Range("C35") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell $ws.Range("D35") "6.90"
Set-TextCell $ws.Range("E35") "  -1.55%  "

# Row 36
Set-TextCell $ws.Range("D36") "5.02"
Set-TextCell $ws.Range("E36") "  +0.16%  "

# Row 37
Set-TextCell $ws.Range("D37") "3.424.72"
Set-TextCell $ws.Range("E37") "  -0.13%  "

# Row 38
Set-TextCell $ws.Range("E38") "  -1.13%  "

# Row 39
Set-TextCell $ws.Range("E39") "  -0.93%  "

# Row 40
Set-TextCell $ws.Range("D40") "25.94"
Set-TextCell $ws.Range("E40") "  -5.10%  "

# Row 41
Set-TextCell $ws.Range("E41") "  -0.29%  "

# Row 42
Set-TextCell $ws.Range("E42") "  -0.30%  "

# Row 43
Set-TextCell $ws.Range("E43") "  -1.21%  "

# Row 44
Set-TextCell $ws.Range("E44") "  +2.04%  "

# Row 45
Set-TextCell $ws.Range("D45") "2.480.36"
Set-TextCell $ws.Range("E45") "  -0.42%  "

# Row 46
Set-TextCell $ws.Range("D46") "22.83"
Set-TextCell $ws.Range("E46") "  -0.95%  "

# Row 47
Set-TextCell $ws.Range("D47") "6.67"
Set-TextCell $ws.Range("E47") "  -2.36%  "

# Row 48
Set-TextCell $ws.Range("E48") "  +0.05%  "

# Row 49
Set-TextCell $ws.Range("E49") "  -1.09%  "

# Row 50
Set-TextCell $ws.Range("D50") "2.02"
Set-TextCell $ws.Range("E50") "  -4.22%  "

# Row 51
Set-TextCell $ws.Range("D51") "0.207"
Set-TextCell $ws.Range("E51") "  -1.71%  "
